# Pinout.xlsx update — rewire some connections now that an IMU is being
# tried (see commit message): UART moved onto D9/D8, I2C SDA/SCL moved onto
# D4/D3 (replacing the old "FPGA UART conn" pins), the motor IN1 pin moved
# off D12 onto A0, ENA moved off D10 onto A1, and A3/A4 (previously the
# unused Radar freq/amplitude breakout) now carry IN3/IN4 motor instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D12 no longer wired to IN1 motor
$ws.Range("C2").ClearContents()

# D9 / D8 now carry UART (replacing IN3 motor / IN4 motor)
$ws.Range("C5").Value = "UART"
$ws.Range("C6").Value = "UART"

# D4 / D3 now carry I2C (replacing the old "FPGA UART conn" pins)
$ws.Range("C10").Value = "I2C SDA"
$ws.Range("C13").Value = "I2C SCL"

# D10 no longer wired to ENA
$ws.Range("C17").ClearContents()

# A0 now carries IN1 motor, A1 now carries ENA
$ws.Range("C18").Value = "IN1 motor"
$ws.Range("C19").Value = "ENA"

# A3 / A4 now carry IN3 motor / IN4 motor (replacing the unused radar pins)
$ws.Range("C21").Value = "IN3 motor"
$ws.Range("C22").Value = "IN4 motor"

# Leave the selection where the author left it when they saved
$ws.Range("C2").Select()
